# Insert a new data row at row 50 (pushes existing rows 50-146 down to 51-147)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()

$ws.Cells.Item(50,1).Value = 4
$ws.Cells.Item(50,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(50,3).Value = "Los Lagos"
$ws.Cells.Item(50,4).Value = 44469
$ws.Cells.Item(50,5).Value = 10
$ws.Cells.Item(50,6).Value = "Fruta"
$ws.Cells.Item(50,7).Value = 100102
$ws.Cells.Item(50,8).Value = "Cítricos"
$ws.Cells.Item(50,9).Value = 100102006
$ws.Cells.Item(50,10).Value = "Pomelo"
$ws.Cells.Item(50,11).Value = "Start Ruby"
$ws.Cells.Item(50,12).Value = "Primera"
$ws.Cells.Item(50,13).Value = 80
$ws.Cells.Item(50,14).Value = 12000
$ws.Cells.Item(50,15).Value = 12000
$ws.Cells.Item(50,16).Value = 12000
$ws.Cells.Item(50,17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(50,18).Value = "Región de O'Higgins"
$ws.Cells.Item(50,19).Value = 857
$ws.Cells.Item(50,20).Value = 14
